# Update recomputed TPM-derived NATMI ligand/receptor metrics (Prn-Rpsa, YoungD4)
# Sets the new values for columns G,H,I,J,M,N,O,P,Q,R,S,T across rows 2-11
# per the refreshed "new tpm" script run (commit: "update scripts wuth new tpm").
$ws = $excel.ActiveWorkbook.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.44393
$ws.Range("H2").Value = 1.33179
$ws.Range("I2").Value = 0.5353249597237094
$ws.Range("J2").Value = 0.5353249597237095
$ws.Range("M2").Value = 127.5808283333333
$ws.Range("N2").Value = 382.742485
$ws.Range("O2").Value = 0.1532286472569342
$ws.Range("P2").Value = 0.1532286472569342
$ws.Range("Q2").Value = 56.63695712201666
$ws.Range("R2").Value = 509.73261409815
$ws.Range("S2").Value = 0.08202711942133677
$ws.Range("T2").Value = 0.08202711942133679
# Row 3
$ws.Range("G3").Value = 0.44393
$ws.Range("H3").Value = 1.33179
$ws.Range("I3").Value = 0.5353249597237094
$ws.Range("J3").Value = 0.5353249597237095
$ws.Range("O3").Value = 0.341528235684153
$ws.Range("P3").Value = 0.341528235684153
$ws.Range("Q3").Value = 126.2369693048767
$ws.Range("R3").Value = 1136.13272374389
$ws.Range("S3").Value = 0.1828285890121287
$ws.Range("T3").Value = 0.1828285890121288
# Row 4
$ws.Range("G4").Value = 0.44393
$ws.Range("H4").Value = 1.33179
$ws.Range("I4").Value = 0.5353249597237094
$ws.Range("J4").Value = 0.5353249597237095
$ws.Range("M4").Value = 155.9545746666667
$ws.Range("N4").Value = 467.863724
$ws.Range("O4").Value = 0.1873064223040503
$ws.Range("P4").Value = 0.1873064223040503
$ws.Range("Q4").Value = 69.23291433177333
$ws.Range("R4").Value = 623.0962289859601
$ws.Range("S4").Value = 0.1002698029759078
$ws.Range("T4").Value = 0.1002698029759078
# Row 5
$ws.Range("G5").Value = 0.44393
$ws.Range("H5").Value = 1.33179
$ws.Range("I5").Value = 0.5353249597237094
$ws.Range("J5").Value = 0.5353249597237095
$ws.Range("M5").Value = 132.7811556666666
$ws.Range("N5").Value = 398.343467
$ws.Range("O5").Value = 0.1594744063806953
$ws.Range("P5").Value = 0.1594744063806954
$ws.Range("Q5").Value = 58.94553843510332
$ws.Range("R5").Value = 530.50984591593
$ws.Range("S5").Value = 0.08537063017270821
$ws.Range("T5").Value = 0.08537063017270824
# Row 6
$ws.Range("G6").Value = 0.44393
$ws.Range("H6").Value = 1.33179
$ws.Range("I6").Value = 0.5353249597237094
$ws.Range("J6").Value = 0.5353249597237095
$ws.Range("M6").Value = 131.9384486666667
$ws.Range("N6").Value = 395.815346
$ws.Range("O6").Value = 0.1584622883741672
$ws.Range("P6").Value = 0.1584622883741672
$ws.Range("Q6").Value = 58.57143551659333
$ws.Range("R6").Value = 527.14291964934
$ws.Range("S6").Value = 0.08482881814162788
$ws.Range("T6").Value = 0.0848288181416279
# Row 7
$ws.Range("G7").Value = 0.385342
$ws.Range("I7").Value = 0.4646750402762905
$ws.Range("J7").Value = 0.4646750402762905
$ws.Range("M7").Value = 127.5808283333333
$ws.Range("N7").Value = 382.742485
$ws.Range("O7").Value = 0.1532286472569342
$ws.Range("P7").Value = 0.1532286472569342
$ws.Range("Q7").Value = 49.16225155162333
$ws.Range("R7").Value = 442.46026396461
$ws.Range("S7").Value = 0.07120152783559741
$ws.Range("T7").Value = 0.07120152783559741
# Row 8
$ws.Range("G8").Value = 0.385342
$ws.Range("I8").Value = 0.4646750402762905
$ws.Range("J8").Value = 0.4646750402762905
$ws.Range("O8").Value = 0.341528235684153
$ws.Range("P8").Value = 0.341528235684153
$ws.Range("R8").Value = 986.1907418577658
$ws.Range("S8").Value = 0.1586996466720242
$ws.Range("T8").Value = 0.1586996466720242
# Row 9
$ws.Range("G9").Value = 0.385342
$ws.Range("I9").Value = 0.4646750402762905
$ws.Range("J9").Value = 0.4646750402762905
$ws.Range("M9").Value = 155.9545746666667
$ws.Range("N9").Value = 467.863724
$ws.Range("O9").Value = 0.1873064223040503
$ws.Range("P9").Value = 0.1873064223040503
$ws.Range("Q9").Value = 60.09584771120267
$ws.Range("R9").Value = 540.8626294008241
$ws.Range("S9").Value = 0.08703661932814244
$ws.Range("T9").Value = 0.08703661932814245
# Row 10
$ws.Range("G10").Value = 0.385342
$ws.Range("I10").Value = 0.4646750402762905
$ws.Range("J10").Value = 0.4646750402762905
$ws.Range("M10").Value = 132.7811556666666
$ws.Range("N10").Value = 398.343467
$ws.Range("O10").Value = 0.1594744063806953
$ws.Range("P10").Value = 0.1594744063806954
$ws.Range("Q10").Value = 51.16615608690466
$ws.Range("R10").Value = 460.495404782142
$ws.Range("S10").Value = 0.07410377620798714
$ws.Range("T10").Value = 0.07410377620798715
# Row 11
$ws.Range("G11").Value = 0.385342
$ws.Range("I11").Value = 0.4646750402762905
$ws.Range("J11").Value = 0.4646750402762905
$ws.Range("M11").Value = 131.9384486666667
$ws.Range("N11").Value = 395.815346
$ws.Range("O11").Value = 0.1584622883741672
$ws.Range("P11").Value = 0.1584622883741672
$ws.Range("Q11").Value = 50.84142568611067
$ws.Range("R11").Value = 457.572831174996
$ws.Range("S11").Value = 0.0736334702325393
$ws.Range("T11").Value = 0.0736334702325393
